$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (date moved from 04-01 to 04-02)
$ws.Name = "Through 2022-04-02"

# Update the header label in I1 to match the new "through" date
$ws.Range("I1").Value = "2022 (through 04-02)"

# Update April 2022 carjacking count (row 5) and the Total row (row 14)
$ws.Range("I5").Value = 6
$ws.Range("I14").Value = 439
